$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.548.73"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.31%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.876.70"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.68%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.016"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +1.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.57"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.41%  "

$ws.Range("E6").Value = "  +0.86%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4790"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.44%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3778"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.63%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07382"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.35%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9399"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.58%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.74"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.64%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07856"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.08%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.900.82"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.05%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.454"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.64%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.599"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.03%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.14"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.89%  "

$ws.Range("E17").Value = "  +0.87%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008961"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.38%  "

$ws.Range("E19").Value = "  +0.83%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.94"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.60%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.581.69"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.34%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.143"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.12%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.78"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.08%  "

$ws.Range("E24").Value = "  +2.64%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.66"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.01%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "18.59"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.38%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.020"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.67%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "116.03"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.56%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.019"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.45%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08939"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.16%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.325"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.47%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.216"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.93%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.615"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.84%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7517"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.65%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.691"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.05%  "

$ws.Range("E36").Value = "  +6.58%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.119"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.77%  "

$ws.Range("E38").Value = "  +0.91%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.012"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.63%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5360"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.51%  "

$ws.Range("E41").Value = "  +2.76%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1525"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.69%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.439"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.56%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.66"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.46%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4845"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.98%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.015"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.81%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.663"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.43%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "103.12"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.12%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "67.33"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.99%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06106"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.24%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9012"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.70%  "
